$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "species" column (column I) with its header
$ws.Range("I1").Value = "species"

# Fill in the species value for every data row (rows 2 through 61)
$ws.Range("I2:I61").Value = "A. elegantissima x B. muscatinei"

# Match the saved selection state recorded in the workbook
$ws.Range("N56").Select()
